# add new word and modify AddShop&EditShop api
# -------------------------------------------------------------------
# Adds two new field-dictionary rows (IS_PRINT / PERSON_NUM) to the
# ORDER_HEAD sheet (sheet6.xml), which introduces four new shared
# strings (IS_PRINT, 是否打印, PERSON_NUM, 人数) and re-uses the existing
# "CHAR" / "DECIMAL(18,0)" shared strings for the type column.
# -------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ORDER_HEAD")

# --- Row 11: IS_PRINT / 是否打印 / CHAR -----------------------------
$ws.Range("A11").Value = "IS_PRINT"
$ws.Range("B11").Value = "是否打印"
$ws.Range("C11").Value = "CHAR"

# --- Row 12: PERSON_NUM / 人数 / DECIMAL(18,0) ----------------------
$ws.Range("A12").Value = "PERSON_NUM"
$ws.Range("B12").Value = "人数"
$ws.Range("C12").Value = "DECIMAL(18,0)"

# Copy the formatting used by the existing type-column cells (C9:C10)
# down onto the two new rows so the style matches the rest of the
# table (centered, bordered "type" cell style).
$ws.Range("C9:C10").Copy()
$ws.Range("C11:C12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the sheet's selection to match the saved view state.
[void]$ws.Activate()
[void]$ws.Range("D28").Select()
